$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 100000
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 90000

$ws.Range("C11").Value = "قابل پرداخت"

$ws.Columns("D").ColumnWidth = 18.625
$ws.Columns("E").ColumnWidth = 13.875

$ws.Range("C11").Select()
